# Update "想去人数" (want-to-go count) figures refreshed from the live
# Bilibili data source for the 展览 (Exhibition), 演出 (Performance) and
# 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# 展览 (sheet 1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1376
$ws.Range("F3").Value = 2222
$ws.Range("F7").Value = 696
$ws.Range("F8").Value = 123
$ws.Range("F11").Value = 2562
$ws.Range("F12").Value = 1620
$ws.Range("F15").Value = 265
$ws.Range("F16").Value = 644
$ws.Range("F17").Value = 820
$ws.Range("F18").Value = 102
$ws.Range("F19").Value = 330
$ws.Range("F20").Value = 1095
$ws.Range("F22").Value = 37
$ws.Range("F24").Value = 5523
$ws.Range("F25").Value = 229
$ws.Range("F26").Value = 850
$ws.Range("F27").Value = 107
$ws.Range("F30").Value = 242
$ws.Range("F34").Value = 793
$ws.Range("F36").Value = 63
$ws.Range("F38").Value = 417
$ws.Range("F39").Value = 1151
$ws.Range("F40").Value = 145
$ws.Range("F41").Value = 110
$ws.Range("F42").Value = 189
$ws.Range("F43").Value = 133
$ws.Range("F44").Value = 94

# 演出 (sheet 2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 792
$ws.Range("F6").Value = 13

# 全部类型 (sheet 4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1376
$ws.Range("F4").Value = 2222
$ws.Range("F9").Value = 696
$ws.Range("F10").Value = 123
$ws.Range("F11").Value = 13
$ws.Range("F15").Value = 2562
$ws.Range("F16").Value = 1620
$ws.Range("F19").Value = 265
$ws.Range("F20").Value = 644
$ws.Range("F22").Value = 820
$ws.Range("F23").Value = 103
$ws.Range("F24").Value = 330
$ws.Range("F25").Value = 1095
$ws.Range("F26").Value = 37
$ws.Range("F28").Value = 5523
$ws.Range("F29").Value = 229
$ws.Range("F30").Value = 850
$ws.Range("F31").Value = 107
$ws.Range("F34").Value = 243
$ws.Range("F38").Value = 793
$ws.Range("F39").Value = 63
$ws.Range("F40").Value = 417
$ws.Range("F41").Value = 1151
$ws.Range("F42").Value = 145
$ws.Range("F43").Value = 110
$ws.Range("F44").Value = 189
$ws.Range("F45").Value = 133
$ws.Range("F46").Value = 94
